$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dry_Weather")

# D2 / E2: wrap the existing pattern names in literal double quotes
$ws.Range("D2").Value = '"dry_weather_hourly"'
$ws.Range("E2").Value = '"dry_weather_monthly"'

# G2: the default data set no longer carries a value here - clear it out
$ws.Range("G2").ClearContents()
